$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 — this shifts the existing rows 5..90 down
# to 6..91 and extends the sheet dimension to A1:R91, matching the diff.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly data record.
$ws.Cells.Item(5, 1).Value = 9
$ws.Cells.Item(5, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(5, 3).Value = "Metropolitana"
$ws.Cells.Item(5, 4).Value = 45092
$ws.Cells.Item(5, 5).Value = 13
$ws.Cells.Item(5, 6).Value = 100112035
$ws.Cells.Item(5, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 43
$ws.Cells.Item(5, 11).Value = 17000
$ws.Cells.Item(5, 12).Value = 19000
$ws.Cells.Item(5, 13).Value = 18023
$ws.Cells.Item(5, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(5, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(5, 16).Value = 1202
$ws.Cells.Item(5, 17).Value = 15
$ws.Cells.Item(5, 18).Value = "Hortaliza"
